$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Shop"
$ws.Range("A2").Value = "My Account"
$ws.Range("A3").Value = "Test Cases"
$ws.Range("A4").Value = "AT Site"
$ws.Range("A5").Value = "Demo Site"
$ws.Range("A6").Value = "0 Items₹0.00"
